$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) "localisation" -> "L" + "ocalisation" (two runs) and add a new
#    bullet item "Respect des normes de codage" right after it.
# ------------------------------------------------------------------

function Find-ParagraphByText($doc, $text) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

$target = Find-ParagraphByText $d "localisation"
if ($target -ne $null) {
    # Replace the paragraph's own content (pPr + runs + bookmark) so that
    # "localisation" is produced from two separate runs: "L" and "ocalisation".
    $xmlSplit = "<w:p $wNs>" +
        "<w:pPr>" +
            "<w:pStyle w:val=`"Corpsdetexte`"/>" +
            "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr>" +
            "<w:spacing w:after=`"0`"/>" +
            "<w:ind w:left=`"709`" w:hanging=`"283`"/>" +
            "<w:jc w:val=`"both`"/>" +
        "</w:pPr>" +
        "<w:r><w:t>L</w:t></w:r>" +
        "<w:r><w:t>ocalisation</w:t></w:r>" +
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
        "<w:bookmarkEnd w:id=`"0`"/>" +
        "</w:p>"
    $target.Range.InsertXML($xmlSplit)

    # Re-resolve the paragraph (its text now reads "Localisation") and
    # insert the new bullet item right after it.
    $target = Find-ParagraphByText $d "Localisation"
    if ($target -ne $null) {
        $endPos = $target.Range.End
        $target.Range.InsertParagraphAfter()

        $paras = $d.Paragraphs
        $blankPara = $null
        for ($i = 1; $i -le $paras.Count; $i++) {
            $p = $paras.Item($i)
            if (($p.Range.Start -eq $endPos) -and ($p.Range.End -eq ($endPos + 1))) {
                $blankPara = $p
                break
            }
        }

        if ($blankPara -ne $null) {
            $xmlNew = "<w:p $wNs>" +
                "<w:pPr>" +
                    "<w:pStyle w:val=`"Corpsdetexte`"/>" +
                    "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr>" +
                    "<w:spacing w:after=`"0`"/>" +
                    "<w:ind w:left=`"709`" w:hanging=`"283`"/>" +
                    "<w:jc w:val=`"both`"/>" +
                "</w:pPr>" +
                "<w:r><w:t>Respect des normes de codage</w:t></w:r>" +
                "</w:p>"
            $blankPara.Range.InsertXML($xmlNew)
        }
    }
}

# ------------------------------------------------------------------
# 2) Footer "print time" cached field result: 08:43 -> 11:02
# ------------------------------------------------------------------

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($f = 1; $f -le 3; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $footer.Range.Find.Execute("08:43", $false, $false, $false, $false, $false, $true, 1, $false, "11:02", 2) | Out-Null
        }
    }
}
